# Generate Report for Handoff
#
# The nightly handoff-generation run produced a fresh xliff handoff for the
# two remaining files in each target language (rows 4-7, which previously
# carried the placeholder "low" priority because they hadn't been handed
# off yet). Now that a handoff xliff exists for them, their Priority bumps
# up to "ht" (matching the already-handed-off rows 2-3) and their "Latest
# Handoff Datetime" is stamped with the generation run's timestamp.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# zh-cn: rows 4-7, column E = Priority, column H = Latest Handoff Datetime
foreach ($r in 4..7) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-09-02 04:33:29"
}

# de-de: rows 4-7, column E = Priority, column H = Latest Handoff Datetime
foreach ($r in 4..7) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-09-02 04:33:34"
}
